$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update data rows: rows 3 and 4 now duplicate the values of row 2 ---
$ws.Range("A3").Value = $ws.Range("A2").Value2
$ws.Range("B3").Value = $ws.Range("B2").Value2
$ws.Range("C3").Value = $ws.Range("C2").Value2
$ws.Range("D3").Value = $ws.Range("D2").Value2
$ws.Range("E3").Value = $ws.Range("E2").Value2
$ws.Range("F3").Value = $ws.Range("F2").Value2
$ws.Range("G3").Value = $ws.Range("G2").Value2

$ws.Range("A4").Value = $ws.Range("A2").Value2
$ws.Range("B4").Value = $ws.Range("B2").Value2
$ws.Range("C4").Value = $ws.Range("C2").Value2
$ws.Range("D4").Value = $ws.Range("D2").Value2
$ws.Range("E4").Value = $ws.Range("E2").Value2
$ws.Range("F4").Value = $ws.Range("F2").Value2
$ws.Range("G4").Value = $ws.Range("G2").Value2

# --- Apply a custom "Phone Number" style number format to column C (Business Phone) ---
$ws.Range("C1:C4").NumberFormat = "[<=9999999]###\-####;\(###\)\ ###\-####"

# --- Update the selection stored in the worksheet view ---
$ws.Range("C1:C1048576").Select()

# --- Switch the sheet's print orientation to Portrait ---
$ws.PageSetup.Orientation = 1
